$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 371, shifting the existing rows 371-385 down to 374-388.
$ws.Rows("371:373").Insert()

# New weekly price report rows (market date 2021-09-22 / serial 44461),
# following the same "1a/2a/3a amarillo" pattern as the other weekly blocks.
$newRows = @(
    @{ Row = 371; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44461; E = 4; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102003; J = "Limón"; K = "Sin especificar"; L = "1a amarillo"; M = 750; N = 3300; O = 3500; P = 3400; Q = "`$/malla 16 kilos"; R = "Provincia de Limarí"; S = 212; T = 16 },
    @{ Row = 372; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44461; E = 4; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102003; J = "Limón"; K = "Sin especificar"; L = "2a amarillo"; M = 600; N = 2300; O = 2500; P = 2400; Q = "`$/malla 16 kilos"; R = "Provincia de Limarí"; S = 150; T = 16 },
    @{ Row = 373; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44461; E = 4; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102003; J = "Limón"; K = "Sin especificar"; L = "3a amarillo"; M = 450; N = 1300; O = 1500; P = 1400; Q = "`$/malla 16 kilos"; R = "Provincia de Limarí"; S = 88; T = 16 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
    $ws.Cells.Item($r, 8).Value = $rowData.H
    $ws.Cells.Item($r, 9).Value = $rowData.I
    $ws.Cells.Item($r, 10).Value = $rowData.J
    $ws.Cells.Item($r, 11).Value = $rowData.K
    $ws.Cells.Item($r, 12).Value = $rowData.L
    $ws.Cells.Item($r, 13).Value = $rowData.M
    $ws.Cells.Item($r, 14).Value = $rowData.N
    $ws.Cells.Item($r, 15).Value = $rowData.O
    $ws.Cells.Item($r, 16).Value = $rowData.P
    $ws.Cells.Item($r, 17).Value = $rowData.Q
    $ws.Cells.Item($r, 18).Value = $rowData.R
    $ws.Cells.Item($r, 19).Value = $rowData.S
    $ws.Cells.Item($r, 20).Value = $rowData.T
}
